$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: both 2021 (E) and 2022 (F) columns get an "X" mark
$ws.Range("E7").Value = "X"
$ws.Range("F7").Value = "X"

# Row 8: both 2021 (E) and 2022 (F) columns get an "X" mark
$ws.Range("E8").Value = "X"
$ws.Range("F8").Value = "X"

# Move the active selection to G4, matching the saved cursor position
$ws.Range("G4").Select()
